$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update benchmark numbers for version 2.4.0 (row 26)
$ws.Range("B26").Value = 1021
$ws.Range("F26").Value = 0.94769999999999999
$ws.Range("G26").Value = 353
$ws.Range("H26").Value = 0.27160000000000001
$ws.Range("I26").Value = 353
$ws.Range("J26").Value = 1.508
$ws.Range("K26").Value = 36
$ws.Range("L26").Value = 0.9698
$ws.Range("M26").Value = 80
$ws.Range("N26").Value = 1.1319999999999999
$ws.Range("O26").Value = 80
$ws.Range("P26").Value = 0.64659999999999995

# Add release note
$ws.Range("R26").Value = "Switched to chem props from ctxR"

# Update the active selection/view position to match the author's final state
$ws.Range("F27").Select()
